$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone column R ("2021") into a brand-new column S so every new S-cell
# picks up the exact same direct formatting (number format / font / fill /
# border) its same-row R neighbour already has - then the values below are
# overwritten with the real 2022 figures.
$ws.Columns("R:R").Copy()
$ws.Columns("S:S").Insert(-4161)
$excel.CutCopyMode = $false

# 2022 values (column S) for rows 3..33.
$values = @{
    3  = 2022
    4  = 57.890663775669807
    5  = 50.022742766269019
    6  = 65.935557126462314
    7  = 41.756312988336546
    8  = 37.159844500343013
    9  = 46.265141318977122
    10 = 54.165137501740752
    11 = 47.516379220948068
    12 = 60.737168699398438
    13 = 38.800251455475774
    14 = 28.662467801253705
    15 = 48.983865238282192
    16 = 52.046737970697684
    17 = 42.039438248006412
    18 = 61.864439546842981
    19 = 54.559895023999445
    20 = 48.115700180781531
    21 = 60.917884338341217
    22 = 44.498872817808376
    23 = 40.811480640517637
    24 = 48.122142747774355
    25 = 86.286397363931727
    26 = 68.335423960459067
    27 = 104.46468013142454
    28 = 63.925654994479551
    29 = 57.785410559676791
    30 = 71.340059495655098
    31 = 50.853410128538314
    32 = 51.223021582733814
    33 = 50.500891999451071
}

foreach ($row in $values.Keys) {
    $ws.Range("S" + $row).Value2 = $values[$row]
}

# Match the recorded selection change (T3 -> T4).
$ws.Range("T4").Select()
